# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (dated 2022-05-25, serial 44706) at the top
# of the "Comercializadora del Agro de Limari - Limon" data block, pushing
# the existing rows (previously 557-577) down to 560-580.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 557; this shifts rows 557:577 -> 560:580
$ws.Rows("557:559").Insert()

# Common (unchanged) values shared by every row in this data block
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102003
$categoria   = "Limón"
$variedad    = "Sin especificar"
$origen      = "Provincia de Limarí"
$fecha       = 44706

$newRows = @(
    @{ Row = 557; Calidad = "1a amarillo"; Volumen = 750; PMin = 4800; PMax = 5000; PProm = 4900; Unidad = "`$/malla 16 kilos"; PrecioKg = 306; KgUnidad = 16 },
    @{ Row = 558; Calidad = "2a amarillo"; Volumen = 750; PMin = 3800; PMax = 4000; PProm = 3900; Unidad = "`$/malla 16 kilos"; PrecioKg = 244; KgUnidad = 16 },
    @{ Row = 559; Calidad = "3a amarillo"; Volumen = 450; PMin = 2800; PMax = 3000; PProm = 2900; Unidad = "`$/malla 16 kilos"; PrecioKg = 181; KgUnidad = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad

    # Match the date-formatted number format used by the rest of column D
    $ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row + 3, 4).NumberFormat
}
